$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.189389586448669
$ws.Range("B1").Value = 2.226623773574829
$ws.Range("C1").Value = 6.506722450256348
$ws.Range("D1").Value = 2.303606748580933
$ws.Range("E1").Value = 1.190868139266968
